# Automatic map update (2025-08-22 08:55:17)
#
# Rows 3-8 get re-shuffled (same 6 records, new order) and two brand-new
# rows (71, 72) are appended at the bottom of the "NEW" sheet.
#
# Columns I (Attachments), M (Coordenada_X) and N (Coordenada_Y) are
# numeric; every other column (A,B,C,D,E,F,G,H,J,K,L,O,P) is text - even
# when its contents look like a number (e.g. Caso/OT ids) - so those
# columns are forced to Text format before the value is written, which
# keeps ids such as "801645368" or "4862" from being auto-converted to
# numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$textCols = @(1,2,3,4,5,6,7,8,10,11,12,15,16)
$numCols  = @(9,13,14)

function Set-RowValues {
    param(
        [int]$Row,
        [object[]]$Values
    )

    foreach ($c in $textCols) {
        $cell = $ws.Cells.Item($Row, $c)
        $cell.NumberFormat = "@"
        $cell.Value = [string]$Values[$c - 1]
    }

    foreach ($c in $numCols) {
        $ws.Cells.Item($Row, $c).Value = [double]$Values[$c - 1]
    }
}

# ---------------------------------------------------------------------
# Rows 3-8: permutation of the existing 6 records (captured here as
# literal values taken from the sheet before the edit).
# ---------------------------------------------------------------------

Set-RowValues 3 @('3839','10/23/2024','PICO 1511','13','798390296','NEW','Pendiente','Poste inclinado','1','Aplomo','Sin equipos','Poste','-58.465596','-34.53627','Saavedra','Capital Norte')
Set-RowValues 4 @('801645368','12/13/2024','San Blas 1809','11','801645368','NEW','Pendiente','Picada','0','Cambio','Sin equipos','Pasante','-58.467767','-34.604588','Paternal','Capital Norte')
Set-RowValues 5 @('5589','12/31/2023','ARCOS 1520','13','799540526','NEW','Pendiente de Traspaso PROPIO','Picada','0','Cambio','Nodo Teco','Pasante','-58.449125','-34.565958','Colegiales','Capital Norte')
Set-RowValues 6 @('4595','1/15/2025','PAROISSIEN 1806','13','802747617','NEW','Pendiente','Aplomar','1','Aplomo','Sin equipos','Terminal','-58.464172','-34.543845','Saavedra','Capital Norte')
Set-RowValues 7 @('4662','1/21/2025','ALTOLAGUIRRE 2397','12','802823938','NEW','Pendiente','Inclinada','1','Aplomo','Sin equipos','Pasante','-58.490766','-34.576987','Paternal','Capital Norte')
Set-RowValues 8 @('4862','1/23/2025','ARCOS 2263','13','802857379','NEW','Pendiente de Traspaso PROPIO','picada','0','Cambio','Nodo Teco','Pasante','-58.455082','-34.558883','Saavedra','Capital Norte')

# ---------------------------------------------------------------------
# New rows 71 and 72 appended at the bottom.
# ---------------------------------------------------------------------

Set-RowValues 71 @('-558','8/21/2025','Blanco Encalada 4210','12','Pendiente ADM','NEW','Pendiente','Colocar columna R400 para pedir taspaso de fuente telecom','1','Cambio','Fuente Teco','Pasante','-58.477593','-34.570321','Colegiales','Capital Norte')
Set-RowValues 72 @('-559','8/21/2025','Av. Del Libertador 6736','13','809098713','NEW','Pendiente','Picada','1','Cambio','Sin equipos','Pasante','-58.453398','-34.550238','Saavedra','Capital Norte')

Write-Output "Updated rows 3-8 and appended rows 71-72 on sheet '$($ws.Name)'."
